# Update the footer page-number text boxes from "n/17" to the corrected
# "n/16" numbering (the deck went from 17 to 16 total content slides).
#
# Mapping of Slide index -> Shape index (within that slide's Shapes
# collection) for the shape holding the page-number text, and the new
# text to assign.

$p = $ppt.ActivePresentation

$pageNumberUpdates = @(
    @{ Slide = 2;  Shape = 6; Text = "2/16" },
    @{ Slide = 3;  Shape = 6; Text = "3/16" },
    @{ Slide = 4;  Shape = 7; Text = "4/16" },
    @{ Slide = 5;  Shape = 7; Text = "5/16" },
    @{ Slide = 6;  Shape = 6; Text = "6/16" },
    @{ Slide = 7;  Shape = 7; Text = "7/16" },
    @{ Slide = 8;  Shape = 5; Text = "8/16" },
    @{ Slide = 9;  Shape = 4; Text = "9/16" },
    @{ Slide = 10; Shape = 3; Text = "10/16" },
    @{ Slide = 11; Shape = 4; Text = "11/16" },
    @{ Slide = 12; Shape = 6; Text = "12/16" },
    @{ Slide = 13; Shape = 5; Text = "13/16" },
    @{ Slide = 14; Shape = 6; Text = "14/16" },
    @{ Slide = 15; Shape = 4; Text = "15/16" },
    @{ Slide = 16; Shape = 4; Text = "16/16" }
)

foreach ($update in $pageNumberUpdates) {
    $slide = $p.Slides.Item($update.Slide)
    $shape = $slide.Shapes.Item($update.Shape)
    $shape.TextFrame.TextRange.Text = $update.Text
}
